# Change design to use CAPD batch update-service job. Generate update-service.csv
#
# Renames the table/header columns from the old *-Status naming to the new
# *-state naming (and reorders the "new" columns), widens column F slightly,
# and moves the active cell selection to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header cells (this also renames the backing ListObject/table
# column names and the shared-string table entries). Set them in reverse
# (F,E,D,C) order so the newly introduced shared strings land in the same
# order as the source workbook.
$ws.Range("F1").Value = " New-EmailVerified-state"
$ws.Range("E1").Value = "NewAccount-state"
$ws.Range("D1").Value = " SecondChance-state"
$ws.Range("C1").Value = " PP-state"

# Column F (New-EmailVerified-state) becomes its own, wider column instead
# of sharing a merged width definition with column E.
$ws.Columns.Item(6).ColumnWidth = 28

# Move the selection/active cell to C2, as last left by the editor.
$ws.Range("C2").Select()
